$wb = $excel.ActiveWorkbook

# --- Sheet: Escapement ---
$ws = $wb.Worksheets.Item("Escapement")
$ws.Cells.Item(60, 1).Value = 878968.868545003
$ws.Cells.Item(60, 2).Value = 4410596.299419
$ws.Cells.Item(60, 3).Value = 4696645.94987

# --- Sheet: Total Catch ---
$ws = $wb.Worksheets.Item("Total Catch")
$ws.Cells.Item(59, 1).Value = 888735.007361009
$ws.Cells.Item(59, 2).Value = 5761028.94391506
$ws.Cells.Item(59, 3).Value = 2397457.90028316
$ws.Cells.Item(60, 1).Value = 1143714.92892355
$ws.Cells.Item(60, 2).Value = 8257698.23101683
$ws.Cells.Item(60, 3).Value = 10373442.6328714

# --- Sheet: Run Size ---
$ws = $wb.Worksheets.Item("Run Size")
$ws.Cells.Item(59, 1).Value = 1212549.00741291
$ws.Cells.Item(59, 2).Value = 8004914.94311406
$ws.Cells.Item(59, 3).Value = 3625516.90068216
$ws.Cells.Item(60, 1).Value = 2022666.92896855
$ws.Cells.Item(60, 2).Value = 12667854.2304358
$ws.Cells.Item(60, 3).Value = 15070741.6327414

# --- Sheet: Run Size no Offshore ---
$ws = $wb.Worksheets.Item("Run Size no Offshore")
$ws.Cells.Item(59, 1).Value = 1207044.6538177
$ws.Cells.Item(59, 2).Value = 7957491.59915901
$ws.Cells.Item(59, 3).Value = 3603507.564969
$ws.Cells.Item(60, 1).Value = 1943192.098233
$ws.Cells.Item(60, 2).Value = 12039608.1678099
$ws.Cells.Item(60, 3).Value = 14289206.5646278
